$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Variável" column (B) for all data rows (2-9) to the new period label
$newLabel = "Diferença 2024/01 - 2023/01"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 2).Value = $newLabel
}

# Row 2: Mato Grosso (name unchanged), new value
$ws.Range("A2").Value = "Mato Grosso"
$ws.Range("C2").Value = 2.58776458440888
$ws.Range("D2").Value = "1º"

# Row 3: Goiás -> Rondônia
$ws.Range("A3").Value = "Rondônia"
$ws.Range("C3").Value = 2.397153944904083
$ws.Range("D3").Value = "2º"

# Row 4: Roraima -> Sergipe
$ws.Range("A4").Value = "Sergipe"
$ws.Range("C4").Value = 1.967690121123582
$ws.Range("D4").Value = "3º"

# Row 5: Paraíba -> Santa Catarina
$ws.Range("A5").Value = "Santa Catarina"
$ws.Range("C5").Value = 1.579398232783426
$ws.Range("D5").Value = "4º"

# Row 6: Minas Gerais -> Paraíba
$ws.Range("A6").Value = "Paraíba"
$ws.Range("C6").Value = 1.299777734113761
$ws.Range("D6").Value = "5º"

# Row 7: São Paulo -> Paraná
$ws.Range("A7").Value = "Paraná"
$ws.Range("C7").Value = 1.228725218594363
$ws.Range("D7").Value = "6º"

# Row 8: Sergipe -> Nordeste, drop the D8 placement value
$ws.Range("A8").Value = "Nordeste"
$ws.Range("C8").Value = 0.3163105243272142
$ws.Range("D8").ClearContents()

# Row 9: Nordeste -> Brasil, new value (used to have no D9, still none)
$ws.Range("A9").Value = "Brasil"
$ws.Range("C9").Value = 0.3953181095616998

# Row 10 (previously Brasil) no longer exists - clear it entirely
$ws.Range("A10:D10").ClearContents()
